{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph right before it) that followed the\n// bibliography entry ending in \"... Thomson Pioneira (2008).\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nconst items = paragraphs.items;\nlet jupiterIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === jupiterText) {\n    jupiterIndex = i;\n    break;\n  }\n}\n\nif (jupiterIndex === -1) {\n  throw new Error(\"Could not locate the 'Ver no Jupiter...' paragraph.\");\n}\n\n// The paragraph right before it should be an empty one that is also part\n// of the block being removed; the one right after it should be the\n// copyright paragraph.\nconst toDelete = [];\nif (jupiterIndex - 1 >= 0 && items[jupiterIndex - 1].text === \"\") {\n  toDelete.push(items[jupiterIndex - 1]);\n}\ntoDelete.push(items[jupiterIndex]);\nif (jupiterIndex + 1 < items.length && items[jupiterIndex + 1].text === copyrightText) {\n  toDelete.push(items[jupiterIndex + 1]);\n}\n\n// Delete from last to first so earlier indices stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph right before it) that followed the\n# bibliography entry ending in \"... Thomson Pioneira (2008).\"\n$d = $word.ActiveDocument\n\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$paras = $d.Paragraphs\n$count = $paras.Count\n$jupiterIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = $paras.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq $jupiterText) {\n        $jupiterIndex = $i\n        break\n    }\n}\n\nif ($jupiterIndex -eq -1) {\n    throw \"Could not locate the 'Ver no Jupiter...' paragraph.\"\n}\n\n$hasLeadingBlank = $false\nif ($jupiterIndex -gt 1) {\n    $prevTxt = $paras.Item($jupiterIndex - 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($prevTxt -eq \"\") {\n        $hasLeadingBlank = $true\n    }\n}\n\n$hasTrailingCopyright = $false\nif ($jupiterIndex + 1 -le $count) {\n    $nextTxt = $paras.Item($jupiterIndex + 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($nextTxt -eq $copyrightText) {\n        $hasTrailingCopyright = $true\n    }\n}\n\n$startIndex = $jupiterIndex\nif ($hasLeadingBlank) {\n    $startIndex = $jupiterIndex - 1\n}\n$endIndex = $jupiterIndex\nif ($hasTrailingCopyright) {\n    $endIndex = $jupiterIndex + 1\n}\n\n# Build one contiguous Range spanning every paragraph to remove (including\n# each paragraph mark) and delete it in a single operation.\n$delRange = $d.Range($paras.Item($startIndex).Range.Start, $paras.Item($endIndex).Range.End)\n$delRange.Delete()\n"}
